# Update quizvragen via Admin
#
# 1) "Wiskunde 3" sheet: a new multiple-choice question ("q2", about
#    cos(alpha)) was inserted as row 3, pushing the former rows 3..20
#    (ids q3..q20) down to rows 4..21 unchanged.
# 2) "AC" sheet: a new test/debug "input" question was appended as row 3.

$wb = $excel.ActiveWorkbook

# --- Sheet "Wiskunde 3": insert new row 3 (q2 - cos formula question) ---
$wsWisk = $wb.Worksheets.Item("Wiskunde 3")

# Push existing rows 3..20 down to 4..21, leaving a blank row 3 behind.
$wsWisk.Rows("3:3").Insert()

$wsWisk.Range("A3").Value = "q2"
$wsWisk.Range("B3").Value = "mc"
$wsWisk.Range("C3").Value = "Goniometrie"
$wsWisk.Range("D3").Value = "Wat is de formule voor cos(α)?"
$wsWisk.Range("E3").Value = '["overstaande / schuine", "aanliggende / schuine", "aanliggende / overstaande"]'
$wsWisk.Range("F3").Value = 1
$wsWisk.Range("G3").Value = "cos(α) = aanliggende / schuine"
$wsWisk.Range("I3").Value = "cos(α)=a/h"
$wsWisk.Range("J3").Value = '["cosinus","basisformule"]'
$wsWisk.Range("K3").Value = 1

# --- Sheet "AC": append new row 3 (debug "input" question) ---
$wsAC = $wb.Worksheets.Item("AC")

$wsAC.Range("B3").Value = "input"
$wsAC.Range("D3").Value = "Werkt het toevoegen nu?"
# Leading apostrophe forces this numeric-looking answer to be stored as text.
$wsAC.Range("F3").Value = "'125"
